$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.312.54"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "2.352.36"
$ws.Range("E3").Value = "  -4.19%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "475.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.47"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +16.30%  "
$ws.Range("D9").Value = "2.352.78"
$ws.Range("E9").Value = "  -4.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0956"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.43"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.324"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.124"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "2.752.40"
$ws.Range("E14").Value = "  -4.39%  "
$ws.Range("D15").Value = "55.229.59"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.96"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.41%  "
$ws.Range("E17").Value = "  -4.07%  "
$ws.Range("D18").Value = "2.355.14"
$ws.Range("E18").Value = "  -4.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "315.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.54"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.40%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.70"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "56.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.394"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.153"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.51%  "
$ws.Range("D28").Value = "2.443.71"
$ws.Range("E28").Value = "  -4.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.10"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.88%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "0.0₃0739"
$ws.Range("E31").Value = "  -4.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.15"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.12"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.09"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.82%  "
$ws.Range("E36").Value = "  -5.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.54"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.807"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -5.90%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.52"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0993"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.37"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.574"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.44%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0516"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.38%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.15"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "249.80"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0221"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.38"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.24%  "
$ws.Range("D50").Value = "1.808.14"
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.94%  "
